$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: sale total for product in row 2 (Dulce - Hershey)
$ws.Range("C10").Formula = "=C2*F2"

# Rows 11-16: sale totals for products in rows 3-8
$ws.Range("C11").Formula = "=C3*F3"
$ws.Range("C12").Formula = "=C4*F4"
$ws.Range("C13").Formula = "=C5*F5"
$ws.Range("C14").Formula = "=C6*F6"
$ws.Range("C15").Formula = "=C7*F7"
$ws.Range("C16").Formula = "=C8*F8"

# Labels, written in the same order the original author typed them
# (total: first, then snack:, bebida:, dulce:) so the shared-string
# table indices line up.
$ws.Range("B22").Value = "total:"
$ws.Range("B18").Value = "snack:"
$ws.Range("B19").Value = "bebida:"
$ws.Range("B20").Value = "dulce:"

# Row 18: snack subtotal
$ws.Range("C18").Formula = "=C11+C16"

# Row 19: bebida subtotal
$ws.Range("C19").Formula = "=C12+C14"

# Row 20: dulce subtotal (after 20% discount) and the discount amount
$ws.Range("C20").Formula = "=(C10+C13+C15)-((C10+C13+C15)*0.2)"
$ws.Range("D20").Formula = "=(C10+C13+C15)*0.2"

# Row 22: grand total
$ws.Range("C22").Formula = "=SUM(C18+C19+C20+D20)"

# Apply integer number format to the discount row only, last, so the
# format doesn't propagate to the dependent total cell C22.
$ws.Range("C20:D20").NumberFormat = "0"

# Update the view to match the selected state after adding rows
$ws.Range("E17").Select()
